$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain literal
# text (matching the source data's formatting, e.g. "1.00" not "1").
# Force Text format first so Excel does not reinterpret them as numbers.
$textForceCells = @(
    'D5',
    'D6',
    'D8',
    'D12',
    'D14',
    'D16',
    'D19',
    'D20',
    'D22',
    'D24',
    'D25',
    'D28',
    'D29',
    'D31',
    'D34',
    'D35',
    'D37',
    'D40',
    'D42',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51',
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '61.455.19'
$ws.Range('E2').Value = '  -1.88%  '
$ws.Range('D3').Value = '2.450.66'
$ws.Range('E3').Value = '  -4.09%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '547.65'
$ws.Range('E5').Value = '  -2.74%  '
$ws.Range('D6').Value = '146.74'
$ws.Range('E6').Value = '  -3.34%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '0.586'
$ws.Range('E8').Value = '  -4.31%  '
$ws.Range('D9').Value = '2.448.74'
$ws.Range('E9').Value = '  -4.28%  '
$ws.Range('E10').Value = '  -5.32%  '
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').Value = '5.41'
$ws.Range('E12').Value = '  -3.04%  '
$ws.Range('E13').Value = '  -5.13%  '
$ws.Range('D14').Value = '26.08'
$ws.Range('E14').Value = '  -4.37%  '
$ws.Range('D15').Value = '2.891.86'
$ws.Range('E15').Value = '  -4.34%  '
$ws.Range('D16').Value = '0.0000168'
$ws.Range('E16').Value = '  -3.85%  '
$ws.Range('D17').Value = '61.420.78'
$ws.Range('E17').Value = '  -1.77%  '
$ws.Range('D18').Value = '2.450.32'
$ws.Range('E18').Value = '  -4.98%  '
$ws.Range('D19').Value = '10.94'
$ws.Range('E19').Value = '  -6.17%  '
$ws.Range('D20').Value = '6.97'
$ws.Range('E20').Value = '  -4.21%  '
$ws.Range('E21').Value = '  -4.45%  '
$ws.Range('D22').Value = '319.20'
$ws.Range('E22').Value = '  -3.62%  '
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').Value = '1.87'
$ws.Range('E24').Value = '  +3.18%  '
$ws.Range('D25').Value = '63.69'
$ws.Range('E25').Value = '  -4.12%  '
$ws.Range('D26').Value = '0.0₃0983'
$ws.Range('E26').Value = '  -8.08%  '
$ws.Range('D27').Value = '2.568.81'
$ws.Range('E27').Value = '  -4.80%  '
$ws.Range('D28').Value = '539.20'
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E30').Value = '  -5.25%  '
$ws.Range('D31').Value = '7.77'
$ws.Range('E31').Value = '  -1.02%  '
$ws.Range('E32').Value = '  -7.01%  '
$ws.Range('E33').Value = '  -5.84%  '
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').Value = '  -4.82%  '
$ws.Range('D35').Value = '1.59'
$ws.Range('E35').Value = '  -3.54%  '
$ws.Range('E36').Value = '  -8.34%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('E38').Value = '  -5.58%  '
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('D40').Value = '18.29'
$ws.Range('E40').Value = '  -4.54%  '
$ws.Range('E41').Value = '  -1.78%  '
$ws.Range('D42').Value = '140.32'
$ws.Range('E42').Value = '  -8.04%  '
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').Value = '40.21'
$ws.Range('E44').Value = '  -3.16%  '
$ws.Range('D45').Value = '2.31'
$ws.Range('E45').Value = '  -2.93%  '
$ws.Range('D46').Value = '142.37'
$ws.Range('E46').Value = '  -7.63%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '21.85'
$ws.Range('E47').Value = '  -4.29%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').Value = '3.61'
$ws.Range('E48').Value = '  -3.75%  '
$ws.Range('D49').Value = '0.0535'
$ws.Range('E49').Value = '  -5.20%  '
$ws.Range('D50').Value = '0.588'
$ws.Range('E50').Value = '  -4.41%  '
$ws.Range('D51').Value = '0.0930'
$ws.Range('E51').Value = '  -4.74%  '
